$d = $word.ActiveDocument
Write-Output "Paragraphs count: $($d.Paragraphs.Count)"
